# Update computed loading_percent results for the 380 kV case (Case_1_52)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 17.92153576870311; "D2" = 7.95272158466683; "E2" = 14.02384515229066; "F2" = 40.9586823234696; "G2" = 47.99626456003916; "H2" = 19.06427011945147; "J2" = 11.03101332699818; "K2" = 12.60783656865077; "L2" = 10.0978393148101; "N2" = 22.24426485161215
    "B3" = 17.83787645489113; "D3" = 7.957751643941527; "E3" = 14.04725200182612; "F3" = 40.99544296408094; "G3" = 47.97452590258293; "H3" = 19.10308792969903; "J3" = 11.04750559117189; "K3" = 12.4158431603883; "L3" = 10.07189754898118; "N3" = 22.30958952790314
    "B4" = 17.78990763820669; "D4" = 7.961457417707553; "E4" = 14.06260893708682; "F4" = 41.02653034269415; "G4" = 47.97340045154926; "H4" = 19.13004070740132; "J4" = 11.05818960207918; "K4" = 12.29932382132781; "L4" = 10.05759489407683; "N4" = 22.35160290079289
    "B5" = 17.77123036132359; "D5" = 7.963123220583837; "E5" = 14.06911526436805; "F5" = 41.04133813750347; "G5" = 47.97601531620174; "H5" = 19.14180765445063; "J5" = 11.06268404714015; "K5" = 12.25224605248902; "L5" = 10.05217937427027; "N5" = 22.36920370135328
    "B6" = 17.76818202542675; "D6" = 7.963409240830566; "E6" = 14.07021064523095; "F6" = 41.04392610281723; "G6" = 47.97663510516937; "H6" = 19.14380884486124; "J6" = 11.06343885105633; "K6" = 12.24445506335595; "L6" = 10.05130518519216; "N6" = 22.37215533042301
    "B7" = 17.78965220557328; "D7" = 7.961479252410178; "E7" = 14.06269567777847; "F7" = 41.02672138660396; "G7" = 47.97342327358406; "H7" = 19.13019622932594; "J7" = 11.05824964583815; "K7" = 12.29868719408386; "L7" = 10.05752018135131; "N7" = 22.35183832610022
    "B8" = 17.89199456642469; "D8" = 7.954328098411942; "E8" = 14.0317117244755; "F8" = 40.96958885181999; "G8" = 47.98623307276193; "H8" = 19.07700686125317; "J8" = 11.03658434371705; "K8" = 12.54139080656237; "L8" = 10.08855976186363; "N8" = 22.26639445799383
    "B9" = 18.11890467868234; "D9" = 7.945183203290615; "E9" = 13.97874317855493; "F9" = 40.92518714382243; "G9" = 48.10823452803359; "H9" = 18.99747571056804; "J9" = 10.99850613107555; "K9" = 13.02546221848913; "L9" = 10.16213169896408; "N9" = 22.11388457971955
    "B10" = 18.30052916089852; "D10" = 7.941412470949386; "E10" = 13.94454251924895; "F10" = 40.93384167920379; "G10" = 48.25666866617829; "H10" = 18.95418354998174; "J10" = 10.97319217995929; "K10" = 13.38247868983951; "L10" = 10.22364706790051; "N10" = 22.01092180400775
    "B11" = 18.38614624789535; "D11" = 7.940331484264994; "E11" = 13.93000042783932; "F11" = 40.94673413625112; "G11" = 48.33685981608835; "H11" = 18.93778244880144; "J11" = 10.96224900756895; "K11" = 13.5444527820176; "L11" = 10.25318555477152; "N11" = 21.96603609292938
    "B12" = 18.41897587446395; "D12" = 7.940012851673449; "E12" = 13.92463925553907; "F12" = 40.95290168902579; "G12" = 48.36903506457298; "H12" = 18.93204560127768; "J12" = 10.95818700960966; "K12" = 13.605666244959; "L12" = 10.26458844705926; "N12" = 21.94931848469826
    "B13" = 18.41188763193655; "D13" = 7.940077448542851; "E13" = 13.92578741201146; "F13" = 40.95151627175021; "G13" = 48.36202534276461; "H13" = 18.93326005216971; "J13" = 10.95905819457649; "K13" = 13.59248921796818; "L13" = 10.2621230682522; "N13" = 21.95290650186246
    "B14" = 18.38883910896412; "D14" = 7.940303455060668; "E14" = 13.92955644568466; "F14" = 40.94721579857294; "G14" = 48.33947072880612; "H14" = 18.93730097542981; "J14" = 10.9619131842549; "K14" = 13.54949164449878; "L14" = 10.25411935934822; "N14" = 21.96465512739727
    "B15" = 18.37477371839542; "D15" = 7.940453689017973; "E15" = 13.93188403481466; "F15" = 40.94474895474223; "G15" = 48.32589048423108; "H15" = 18.93983788273889; "J15" = 10.96367260911806; "K15" = 13.52313667376939; "L15" = 10.24924496542565; "N15" = 21.97188788305216
    "B16" = 18.29499225383377; "D16" = 7.941495837723436; "E16" = 13.9455132807042; "F16" = 40.93317915559005; "G16" = 48.25168183438438; "H16" = 18.95532169541816; "J16" = 10.97391882774705; "K16" = 13.37187937058289; "L16" = 10.22174741572535; "N16" = 22.0138943836722
    "B17" = 18.24680032787457; "D17" = 7.942297277117848; "E17" = 13.9541342420781; "F17" = 40.92837380142349; "G17" = 48.20939322540021; "H17" = 18.96566417843442; "J17" = 10.98035087469431; "K17" = 13.27893318844977; "L17" = 10.20527241753216; "N17" = 22.04016329553357
    "B18" = 18.2193649172132; "D18" = 7.942817981716391; "E18" = 13.95918844716856; "F18" = 40.92645301604262; "G18" = 48.18626320783502; "H18" = 18.97192280974778; "J18" = 10.98410430342496; "K18" = 13.22543557985948; "L18" = 10.19594317722611; "N18" = 22.05545633144744
    "B19" = 18.21012508894536; "D19" = 7.943004558508966; "E19" = 13.96091615824458; "F19" = 40.9259475350914; "G19" = 48.17863708505478; "H19" = 18.97409508787757; "J19" = 10.98538441543231; "K19" = 13.2073176478506; "L19" = 10.19280984648737; "N19" = 22.06066590251687
    "B20" = 18.25190127925453; "D20" = 7.942205783502529; "E20" = 13.95320662955892; "F20" = 40.92879809299063; "G20" = 48.21377150760836; "H20" = 18.96453112579226; "J20" = 10.97966059790756; "K20" = 13.28883178364419; "L20" = 10.20701106642128; "N20" = 22.03734790525511
    "B21" = 18.39559811981386; "D21" = 7.940234613601116; "E21" = 13.9284454414304; "F21" = 40.94844408784216; "G21" = 48.34604659585532; "H21" = 18.93610119427726; "J21" = 10.96107238340264; "K21" = 13.56212487022234; "L21" = 10.25646439639121; "N21" = 21.96119669105963
    "B22" = 18.49188131493523; "D22" = 7.939474846154032; "E22" = 13.9131110290568; "F22" = 40.96877479835422; "G22" = 48.44303026253734; "H22" = 18.92028291749765; "J22" = 10.94940138348613; "K22" = 13.73999671927868; "L22" = 10.29004853959804; "N22" = 21.91305688566331
    "B23" = 18.44028396972775; "D23" = 7.93983215949824; "E23" = 13.92121782110995; "F23" = 40.95723949984831; "G23" = 48.39030926163675; "H23" = 18.92847257017386; "J23" = 10.95558684109092; "K23" = 13.64515026113748; "L23" = 10.27201058217687; "N23" = 21.93860128767669
    "B24" = 18.24959429400571; "D24" = 7.942246961023305; "E24" = 13.95362569821291; "F24" = 40.92860364815329; "G24" = 48.21178840249792; "H24" = 18.96504240515076; "J24" = 10.97997249895812; "K24" = 13.28435681603493; "L24" = 10.2062245790468; "N24" = 22.03862014909327
    "B25" = 18.05481927349198; "D25" = 7.947137507365648; "E25" = 13.99224205736835; "F25" = 40.92995025683437; "G25" = 48.06487598112221; "H25" = 19.01633468151217; "J25" = 11.00833802635701; "K25" = 12.89401160967998; "L25" = 10.14089637784772; "N25" = 22.15354075746097
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

